$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format first so values such as
# "217.70", "4.450", "0.06271" round-trip as literal strings instead
# of being auto-parsed into numbers (which would also strip trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.957.99"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.644.49"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "217.70"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "0.5243"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "0.2619"
$ws.Range("E8").Value = "  -2.07%  "
$ws.Range("D9").Value = "0.06271"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").Value = "20.38"
$ws.Range("E10").Value = "  -3.44%  "
$ws.Range("D11").Value = "0.07743"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.655.20"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.450"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "0.5440"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "0.0₅8066"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").Value = "64.69"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "25.983.82"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "4.548"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("D20").Value = "191.75"
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "5.972"
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").Value = "139.73"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("D25").Value = "0.1238"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "7.255"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "16.17"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "1.419"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").Value = "0.05937"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").Value = "1.272"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("D31").Value = "3.486"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").Value = "3.234"
$ws.Range("E32").Value = "  -3.26%  "
$ws.Range("D33").Value = "1.526"
$ws.Range("E33").Value = "  -7.74%  "
$ws.Range("D34").Value = "2.413"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "0.9401"
$ws.Range("E35").Value = "  -4.17%  "
$ws.Range("D36").Value = "2.742"
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("D37").Value = "0.5727"
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("D39").Value = "5.851"
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("D40").Value = "0.8463"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "100.54"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").Value = "1.001.84"
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("D44").Value = "1.785.11"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").Value = "56.55"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "0.4285"
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("D49").Value = "1.475"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("D51").Value = "7.797"
$ws.Range("E51").Value = "  -3.70%  "
